$wb = $excel.ActiveWorkbook

# SF User list changes - 16 May - Initial
# Replace the user on the Users sheet: "Ashley Choi" -> "Aadarsh Patel"
$users = $wb.Worksheets.Item("Users")
$users.Range("A2").Value = "Aadarsh Patel"

# Bring the Users sheet to the front / make it the active tab, and move
# the selection to B12 (matches the saved view state in the workbook)
$users.Activate()
$users.Range("B12").Select()
